# TimeSheet.xlsx update: add "What I worked on" column with notes.
$wb = $excel.ActiveWorkbook

# --- Sheet "Provencher" (sheet1): column D ---
$ws1 = $wb.Worksheets.Item("Provencher")
$ws1.Range("D1").Value = "What I worked on"
$ws1.Range("D2").Value = "Tutorials + implementation for Random map Gen"
$ws1.Range("D3").Value = "Implementation of Random Map Gen + Player programming"
$ws1.Columns.Item(4).ColumnWidth = 55
$ws1.Range("D1").Select()

# --- Sheet "Philippona" (sheet2): column D header only ---
$ws2 = $wb.Worksheets.Item("Philippona")
$ws2.Range("D1").Value = "What I worked on"
$ws2.Range("D1").Select()

# --- Sheet "Ningge" (sheet3): column C header only ---
$ws3 = $wb.Worksheets.Item("Ningge")
$ws3.Range("C1").Value = "What I worked on"
$ws3.Range("C1").Select()

$ws1.Select()
